$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row 2 ("Check AI eat again" row) ---
# Cell 3 ("In Progress" column) currently holds "NEED to Test" -> clear it out.
$cell3 = $t.Rows.Item(2).Cells.Item(3)
$cell3.Range.Text = ""

# Cell 4 ("Done" column) currently empty -> set it to "yes".
$cell4 = $t.Rows.Item(2).Cells.Item(4)
$cell4.Range.Text = "yes"

# --- Move the hidden _GoBack bookmark ---
# It currently lives inside Row 4 Cell 1 ("Game Finish? ..." paragraph),
# right before the closing ")". Move it into the empty "In Progress" cell
# of Row 3 ("Make Animation to all move and eat function" row).
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$targetCell = $t.Rows.Item(3).Cells.Item(3)
$targetRange = $targetCell.Range
$targetRange.Collapse(0)
$targetRange.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $targetRange)
